$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '47.445.95'
$ws.Range("E2").Value = '  +2.86%  '
$ws.Range("D3").Value = '2.509.86'
$ws.Range("E3").Value = '  +2.30%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.39'
$ws.Range("E5").Value = '  +0.66%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '109.60'
$ws.Range("E6").Value = '  +4.25%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.528'
$ws.Range("E7").Value = '  +1.62%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  +1.21%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.18'
$ws.Range("E10").Value = '  +8.38%  '
$ws.Range("E11").Value = '  +1.53%  '
$ws.Range("E12").Value = '  +1.05%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.68'
$ws.Range("E13").Value = '  +1.76%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.23'
$ws.Range("E14").Value = '  +2.30%  '
$ws.Range("D15").Value = '2.901.57'
$ws.Range("E15").Value = '  +2.11%  '
$ws.Range("D16").Value = '2.507.33'
$ws.Range("E16").Value = '  +1.71%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.867'
$ws.Range("E17").Value = '  +2.78%  '
$ws.Range("D18").Value = '47.397.75'
$ws.Range("E18").Value = '  +2.95%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.98'
$ws.Range("E19").Value = '  +3.12%  '
$ws.Range("E20").Value = '  +4.39%  '
$ws.Range("D21").Value = '0.0₃0950'
$ws.Range("E21").Value = '  +1.32%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.62'
$ws.Range("E22").Value = '  +11.28%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.02'
$ws.Range("E23").Value = '  -0.90%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '250.10'
$ws.Range("E24").Value = '  +0.80%  '
$ws.Range("E25").Value = '  +3.65%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.24'
$ws.Range("E26").Value = '  +0.66%  '
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("E28").Value = '  +1.61%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.08'
$ws.Range("E29").Value = '  +3.92%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.02'
$ws.Range("E30").Value = '  +6.20%  '
$ws.Range("E31").Value = '  +4.49%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '50.29'
$ws.Range("E32").Value = '  +2.11%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.00'
$ws.Range("E33").Value = '  -1.51%  '
$ws.Range("E34").Value = '  +2.69%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0796'
$ws.Range("E35").Value = '  +4.31%  '
$ws.Range("E36").Value = '  +0.11%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.01'
$ws.Range("E37").Value = '  +5.66%  '
$ws.Range("E38").Value = '  +4.48%  '
$ws.Range("E39").Value = '  +2.90%  '
$ws.Range("E40").Value = '  +1.71%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '123.61'
$ws.Range("E41").Value = '  -3.44%  '
$ws.Range("E42").Value = '  -1.38%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.51'
$ws.Range("E43").Value = '  +2.61%  '
$ws.Range("E44").Value = '  +2.34%  '
$ws.Range("D45").Value = '1.997.96'
$ws.Range("E45").Value = '  +1.94%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.10'
$ws.Range("E46").Value = '  +4.12%  '
$ws.Range("E47").Value = '  -1.43%  '
$ws.Range("E48").Value = '  -2.03%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.09'
$ws.Range("E49").Value = '  -1.14%  '
$ws.Range("E50").Value = '  +8.18%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '78.63'
$ws.Range("E51").Value = '  +1.31%  '
